{"js": "// Replace the worksheet date and the multiplication problems with the new\n// values from the next day's worksheet. Every \"old\" string below is unique\n// within the document, so a plain body-wide search/replace is unambiguous.\nconst replacements = [\n  [\"2024-01-16 Tuesday\", \"2024-01-17 Wednesday\"],\n  [\"72\u00d769=\", \"51\u00d778=\"],\n  [\"65\u00d742=\", \"72\u00d788=\"],\n  [\"70\u00d744=\", \"88\u00d742=\"],\n  [\"68\u00d766=\", \"73\u00d740=\"],\n  [\"61\u00d750=\", \"48\u00d760=\"],\n  [\"11\u00d757=\", \"11\u00d795=\"],\n  [\"78\u00d730=\", \"33\u00d724=\"],\n  [\"51\u00d718=\", \"59\u00d756=\"],\n  [\"35\u00d782=\", \"35\u00d784=\"],\n  [\"76\u00d737=\", \"41\u00d716=\"],\n  [\"18\u00d786=\", \"63\u00d772=\"],\n  [\"97\u00d721=\", \"35\u00d797=\"],\n  [\"52\u00d713=\", \"88\u00d732=\"],\n  [\"63\u00d719=\", \"47\u00d797=\"],\n  [\"34\u00d743=\", \"60\u00d734=\"],\n  [\"94\u00d772=\", \"84\u00d747=\"],\n  [\"67\u00d757=\", \"44\u00d716=\"],\n  [\"56\u00d758=\", \"86\u00d730=\"],\n  [\"29\u00d780=\", \"69\u00d712=\"],\n  [\"28\u00d732=\", \"53\u00d718=\"],\n  [\"50\u00d737=\", \"81\u00d796=\"],\n  [\"88\u00d750=\", \"25\u00d776=\"],\n  [\"40\u00d762=\", \"65\u00d758=\"],\n  [\"75\u00d770=\", \"85\u00d776=\"],\n  [\"80\u00d769=\", \"40\u00d782=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the worksheet date and the multiplication problems with the new\n# values from the next day's worksheet. Every \"old\" string below is unique\n# within the document, so a plain document-wide Find/Replace is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-01-16 Tuesday\", \"2024-01-17 Wednesday\"),\n    @(\"72\u00d769=\", \"51\u00d778=\"),\n    @(\"65\u00d742=\", \"72\u00d788=\"),\n    @(\"70\u00d744=\", \"88\u00d742=\"),\n    @(\"68\u00d766=\", \"73\u00d740=\"),\n    @(\"61\u00d750=\", \"48\u00d760=\"),\n    @(\"11\u00d757=\", \"11\u00d795=\"),\n    @(\"78\u00d730=\", \"33\u00d724=\"),\n    @(\"51\u00d718=\", \"59\u00d756=\"),\n    @(\"35\u00d782=\", \"35\u00d784=\"),\n    @(\"76\u00d737=\", \"41\u00d716=\"),\n    @(\"18\u00d786=\", \"63\u00d772=\"),\n    @(\"97\u00d721=\", \"35\u00d797=\"),\n    @(\"52\u00d713=\", \"88\u00d732=\"),\n    @(\"63\u00d719=\", \"47\u00d797=\"),\n    @(\"34\u00d743=\", \"60\u00d734=\"),\n    @(\"94\u00d772=\", \"84\u00d747=\"),\n    @(\"67\u00d757=\", \"44\u00d716=\"),\n    @(\"56\u00d758=\", \"86\u00d730=\"),\n    @(\"29\u00d780=\", \"69\u00d712=\"),\n    @(\"28\u00d732=\", \"53\u00d718=\"),\n    @(\"50\u00d737=\", \"81\u00d796=\"),\n    @(\"88\u00d750=\", \"25\u00d776=\"),\n    @(\"40\u00d762=\", \"65\u00d758=\"),\n    @(\"75\u00d770=\", \"85\u00d776=\"),\n    @(\"80\u00d769=\", \"40\u00d782=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
